# Apply weekly fruit/vegetable price update: the Date/Calidad/Volumen/Prices/
# Unidad/Origen/Precio-Kg/Kg-unidad columns (D, L-T) for rows 2-8 are re-shuffled
# to reflect the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2 = @{ D = 44742; L = 'Segunda'; M = 100; N = 14000; O = 15000; P = 14500; Q = '$/caja 18 kilos granel';    R = "Región de O'Higgins"; S = 806;   T = 18 }
    3 = @{ D = 44334; L = 'Primera'; M = 100; N = 11000; O = 12000; P = 11500; Q = '$/caja 12 kilos granel';    R = "Región de O'Higgins"; S = 11500; T = 1  }
    4 = @{ D = 44708; L = 'Primera'; M = 70;  N = 12000; O = 13000; P = 12571; Q = '$/caja 12 kilos empedrada'; R = 'Provincia de Curicó';  S = 1048;  T = 12 }
    5 = @{ D = 44714; L = 'Primera'; M = 100; N = 14000; O = 15000; P = 14500; Q = '$/caja 18 kilos granel';    R = "Región de O'Higgins"; S = 806;   T = 18 }
    6 = @{ D = 44330; L = 'Primera'; M = 100; N = 15000; O = 16000; P = 15500; Q = '$/caja 18 kilos granel';    R = 'Provincia de Curicó';  S = 861;   T = 18 }
    7 = @{ D = 44707; L = 'Primera'; M = 60;  N = 12000; O = 13000; P = 12500; Q = '$/caja 12 kilos empedrada'; R = 'Provincia de Curicó';  S = 1042;  T = 12 }
    8 = @{ D = 44719; L = 'Primera'; M = 50;  N = 14000; O = 15000; P = 14400; Q = '$/caja 18 kilos granel';    R = 'Región del Maule';     S = 800;   T = 18 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("D$r").Value = $data.D
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("O$r").Value = $data.O
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
    $ws.Range("S$r").Value = $data.S
    $ws.Range("T$r").Value = $data.T
}
